# Pengerjaan control dan module capaian vaksinasi
# Re-order the Kab/Kota list in column B of the "Data" sheet so it follows
# the vaccination-coverage module ordering (cities first, then the
# remaining regencies), and move the active selection from D11 to D4,
# scrolling the window so row 2 becomes the top visible row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New Kab/Kota order for B3:B21 (rows 3-21)
$kabKota = @(
    "13.71 - KOTA PADANG",
    "13.74 - KOTA PADANG PANJANG",
    "13.75 - KOTA BUKITTINGGI",
    "13.76 - KOTA PAYAKUMBUH",
    "13.72 - KOTA SOLOK",
    "13.73 - KOTA SAWAHLUNTO",
    "13.08 - PASAMAN",
    "13.05 - PADANG PARIAMAN",
    "13.06 - AGAM",
    "13.07 - LIMA PULUH KOTA",
    "13.02 - SOLOK",
    "13.04 - TANAH DATAR",
    "13.03 - SIJUNJUNG",
    "13.01 - PESISIR SELATAN",
    "13.09 - KEPULAUAN MENTAWAI",
    "13.77 - KOTA PARIAMAN",
    "13.12 - PASAMAN BARAT",
    "13.10 - DHARMASRAYA",
    "13.11 - SOLOK SELATAN"
)

$startRow = 3
for ($i = 0; $i -lt $kabKota.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $kabKota[$i]
}

# Move the selection to D4 and scroll the sheet view up so row 2 is the
# top-left visible row (was previously scrolled further down to show D11).
$ws.Activate()
$ws.Range("D4").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
